$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update year header labels (shift 1396-1400 range to 1397-1401) ---
$ws.Range("E8").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E27").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F27").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G27").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H27").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I27").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E34").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F34").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G34").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H34").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I34").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E41").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F41").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G41").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H41").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I41").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E48").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F48").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G48").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H48").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I48").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E55").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F55").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G55").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H55").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I55").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E62").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F62").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G62").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H62").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I62").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E69").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F69").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G69").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H69").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I69").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E76").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F76").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G76").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H76").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I76").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E83").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F83").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G83").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H83").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I83").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E89").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F89").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G89").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H89").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I89").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E95").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F95").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G95").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H95").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I95").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E101").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F101").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G101").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H101").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I101").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("E107").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F107").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G107").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H107").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I107").Value = "دوازده ماهه منتهی به 1401/12"

# --- Shift yearly data columns left and append new year (1401) figures ---
$ws.Range("E10").Value = 576827
$ws.Range("F10").Value = 859182
$ws.Range("G10").Value = 968891
$ws.Range("H10").Value = 2035759
$ws.Range("I10").Value = 6009670
$ws.Range("E11").Value = 53885
$ws.Range("F11").Value = 81403
$ws.Range("G11").Value = 107185
$ws.Range("H11").Value = 156082
$ws.Range("I11").Value = 239562
$ws.Range("E12").Value = 143577
$ws.Range("F12").Value = 203721
$ws.Range("G12").Value = 324922
$ws.Range("H12").Value = 590456
$ws.Range("I12").Value = 829807
$ws.Range("E13").Value = 774289
$ws.Range("F13").Value = 1144306
$ws.Range("G13").Value = 1400998
$ws.Range("H13").Value = 2782297
$ws.Range("I13").Value = 7079039
$ws.Range("E15").Value = 774289
$ws.Range("F15").Value = 1144306
$ws.Range("G15").Value = 1400998
$ws.Range("H15").Value = 2782297
$ws.Range("I15").Value = 7079039
$ws.Range("E16").Value = 20389
$ws.Range("F16").Value = -41442
$ws.Range("G16").Value = -133155
$ws.Range("H16").Value = 10039
$ws.Range("I16").Value = -134911
$ws.Range("E18").Value = 794678
$ws.Range("F18").Value = 1102864
$ws.Range("G18").Value = 1267843
$ws.Range("H18").Value = 2792336
$ws.Range("I18").Value = 6944128
$ws.Range("E19").Value = 39044
$ws.Range("F19").Value = 39619
$ws.Range("G19").Value = 80655
$ws.Range("H19").Value = 148132
$ws.Range("I19").Value = 812828
$ws.Range("E20").Value = -39619
$ws.Range("F20").Value = -80655
$ws.Range("G20").Value = -148132
$ws.Range("H20").Value = -812828
$ws.Range("I20").Value = -1850552
$ws.Range("E21").Value = 794103
$ws.Range("F21").Value = 1061828
$ws.Range("G21").Value = 1200366
$ws.Range("H21").Value = 2127640
$ws.Range("I21").Value = 5906404
$ws.Range("E23").Value = 794103
$ws.Range("F23").Value = 1061828
$ws.Range("G23").Value = 1200366
$ws.Range("H23").Value = 2127640
$ws.Range("I23").Value = 5906404
$ws.Range("E29").Value = 516173
$ws.Range("F29").Value = 521246
$ws.Range("G29").Value = 1102023
$ws.Range("H29").Value = 997225
$ws.Range("I29").Value = 1507312
$ws.Range("E30").Value = 516173
$ws.Range("F30").Value = 521246
$ws.Range("G30").Value = 1102023
$ws.Range("H30").Value = 997225
$ws.Range("I30").Value = 1507312
$ws.Range("E36").Value = 1654944
$ws.Range("F36").Value = 2586309
$ws.Range("G36").Value = 2629690
$ws.Range("H36").Value = 4244705
$ws.Range("I36").Value = 1945968
$ws.Range("E37").Value = 1654944
$ws.Range("F37").Value = 2586309
$ws.Range("G37").Value = 2629690
$ws.Range("H37").Value = 4244705
$ws.Range("I37").Value = 1945968
$ws.Range("E43").Value = 1649871
$ws.Range("F43").Value = 2005532
$ws.Range("G43").Value = 2734488
$ws.Range("H43").Value = 3734618
$ws.Range("I43").Value = 2400427
$ws.Range("E44").Value = 1649871
$ws.Range("F44").Value = 2005532
$ws.Range("G44").Value = 2734488
$ws.Range("H44").Value = 3734618
$ws.Range("I44").Value = 2400427
$ws.Range("E50").Value = 521246
$ws.Range("F50").Value = 1102023
$ws.Range("G50").Value = 997225
$ws.Range("H50").Value = 1507312
$ws.Range("I50").Value = 1052853
$ws.Range("E51").Value = 521246
$ws.Range("F51").Value = 1102023
$ws.Range("G51").Value = 997225
$ws.Range("H51").Value = 1507312
$ws.Range("I51").Value = 1052853
$ws.Range("E57").Value = 64625
$ws.Range("F57").Value = 78652
$ws.Range("G57").Value = 288312
$ws.Range("H57").Value = 501833
$ws.Range("I57").Value = 794193
$ws.Range("E58").Value = 64625
$ws.Range("F58").Value = 78652
$ws.Range("G58").Value = 288312
$ws.Range("H58").Value = 501833
$ws.Range("I58").Value = 794193
$ws.Range("E64").Value = 590854
$ws.Range("F64").Value = 1068842
$ws.Range("G64").Value = 1182412
$ws.Range("H64").Value = 2328119
$ws.Range("I64").Value = 6320907
$ws.Range("E65").Value = 590854
$ws.Range("F65").Value = 1068842
$ws.Range("G65").Value = 1182412
$ws.Range("H65").Value = 2328119
$ws.Range("I65").Value = 6320907
$ws.Range("E71").Value = 576827
$ws.Range("F71").Value = 859182
$ws.Range("G71").Value = 968891
$ws.Range("H71").Value = 2035759
$ws.Range("I71").Value = 6009670
$ws.Range("E72").Value = 576827
$ws.Range("F72").Value = 859182
$ws.Range("G72").Value = 968891
$ws.Range("H72").Value = 2035759
$ws.Range("I72").Value = 6009670
$ws.Range("E78").Value = 78652
$ws.Range("F78").Value = 288312
$ws.Range("G78").Value = 501833
$ws.Range("H78").Value = 794193
$ws.Range("I78").Value = 1105430
$ws.Range("E79").Value = 78652
$ws.Range("F79").Value = 288312
$ws.Range("G79").Value = 501833
$ws.Range("H79").Value = 794193
$ws.Range("I79").Value = 1105430
$ws.Range("E85").Value = 125200
$ws.Range("F85").Value = 150892
$ws.Range("G85").Value = 261621
$ws.Range("H85").Value = 503229
$ws.Range("I85").Value = 526894
$ws.Range("E91").Value = 357024
$ws.Range("F91").Value = 413269
$ws.Range("G91").Value = 449639
$ws.Range("H91").Value = 548476
$ws.Range("I91").Value = 3248207
$ws.Range("E97").Value = 349619
$ws.Range("F97").Value = 428406
$ws.Range("G97").Value = 354323
$ws.Range("H97").Value = 545105
$ws.Range("I97").Value = 2503584
$ws.Range("E103").Value = 150892
$ws.Range("F103").Value = 261621
$ws.Range("G103").Value = 503229
$ws.Range("H103").Value = 526894
$ws.Range("I103").Value = 1049938
$ws.Range("E109").Value = 622
$ws.Range("F109").Value = 805
$ws.Range("G109").Value = 1746
$ws.Range("H109").Value = 5963
$ws.Range("I109").Value = 13222
$ws.Range("E113").Value = 11073
$ws.Range("F113").Value = 16858
$ws.Range("G113").Value = 20680
$ws.Range("H113").Value = 37922
$ws.Range("I113").Value = 67557
$ws.Range("E114").Value = 4142
$ws.Range("F114").Value = 4594
$ws.Range("G114").Value = 4795
$ws.Range("H114").Value = 12036
$ws.Range("I114").Value = 13723
$ws.Range("E115").Value = 12092
$ws.Range("F115").Value = 11497
$ws.Range("G115").Value = 18212
$ws.Range("H115").Value = 30108
$ws.Range("I115").Value = 34173
$ws.Range("E116").Value = 74393
$ws.Range("F116").Value = 114024
$ws.Range("G116").Value = 162436
$ws.Range("H116").Value = 259901
$ws.Range("I116").Value = 374375
$ws.Range("E118").Value = 41255
$ws.Range("F118").Value = 55943
$ws.Range("G118").Value = 117053
$ws.Range("H118").Value = 244526
$ws.Range("I118").Value = 326757
$ws.Range("E119").Value = 143577
$ws.Range("F119").Value = 203721
$ws.Range("G119").Value = 324922
$ws.Range("H119").Value = 590456
$ws.Range("I119").Value = 829807
